$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A79").Value = "2025-04-29 11:50:23"
$ws.Range("B79").Value = 245
